$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44389
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 18
$ws.Range("O2").Value = 200000
$ws.Range("P2").Value = 200000
$ws.Range("R2").Value = 'Provincia de Quillota'
$ws.Range("S2").Value = 571
$ws.Range("D3").Value = 44363
$ws.Range("O3").Value = 230000
$ws.Range("P3").Value = 215000
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 614
$ws.Range("D4").Value = 44196
$ws.Range("K4").Value = 'Red Blush'
$ws.Range("M4").Value = 12
$ws.Range("N4").Value = 130000
$ws.Range("O4").Value = 130000
$ws.Range("P4").Value = 130000
$ws.Range("R4").Value = 'Provincia de Limarí'
$ws.Range("S4").Value = 371
$ws.Range("D5").Value = 44201
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 8
$ws.Range("N5").Value = 200000
$ws.Range("O5").Value = 200000
$ws.Range("P5").Value = 200000
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 571
$ws.Range("D6").Value = 44201
$ws.Range("M6").Value = 16
$ws.Range("N6").Value = 170000
$ws.Range("O6").Value = 170000
$ws.Range("P6").Value = 170000
$ws.Range("S6").Value = 486
$ws.Range("D7").Value = 44208
$ws.Range("N7").Value = 180000
$ws.Range("O7").Value = 180000
$ws.Range("P7").Value = 180000
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("S7").Value = 514
$ws.Range("D8").Value = 44298
$ws.Range("L8").Value = 'Especial'
$ws.Range("M8").Value = 15
$ws.Range("N8").Value = 450000
$ws.Range("O8").Value = 450000
$ws.Range("P8").Value = 450000
$ws.Range("S8").Value = 1286
$ws.Range("D9").Value = 44298
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 430000
$ws.Range("O9").Value = 430000
$ws.Range("P9").Value = 430000
$ws.Range("S9").Value = 1229
$ws.Range("D10").Value = 44446
$ws.Range("M10").Value = 14
$ws.Range("N10").Value = 150000
$ws.Range("O10").Value = 160000
$ws.Range("P10").Value = 155000
$ws.Range("Q10").Value = '$/bins (350 kilos)'
$ws.Range("S10").Value = 443
$ws.Range("T10").Value = 350
$ws.Range("D11").Value = 44195
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 20
$ws.Range("O11").Value = 210000
$ws.Range("P11").Value = 206000
$ws.Range("S11").Value = 589
$ws.Range("D12").Value = 44308
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = 280000
$ws.Range("O12").Value = 280000
$ws.Range("P12").Value = 280000
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 800
$ws.Range("D15").Value = 44356
$ws.Range("M15").Value = 24
$ws.Range("N15").Value = 200000
$ws.Range("O15").Value = 230000
$ws.Range("P15").Value = 215000
$ws.Range("S15").Value = 614
$ws.Range("D16").Value = 44309
$ws.Range("M16").Value = 16
$ws.Range("N16").Value = 350000
$ws.Range("O16").Value = 350000
$ws.Range("P16").Value = 350000
$ws.Range("Q16").Value = '$/bins (350 kilos)'
$ws.Range("R16").Value = 'Región Metropolitana'
$ws.Range("S16").Value = 1000
$ws.Range("T16").Value = 350
$ws.Range("D17").Value = 44376
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 180000
$ws.Range("O17").Value = 180000
$ws.Range("P17").Value = 180000
$ws.Range("R17").Value = 'Hijuelas'
$ws.Range("S17").Value = 514
$ws.Range("D18").Value = 44376
$ws.Range("K18").Value = 'Start Ruby'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 16
$ws.Range("N18").Value = 150000
$ws.Range("O18").Value = 150000
$ws.Range("P18").Value = 150000
$ws.Range("S18").Value = 429
$ws.Range("D19").Value = 44400
$ws.Range("M19").Value = 140
$ws.Range("N19").Value = 9800
$ws.Range("O19").Value = 9800
$ws.Range("P19").Value = 9800
$ws.Range("Q19").Value = '$/caja 14 kilos empedrada'
$ws.Range("R19").Value = 'Región de O''Higgins'
$ws.Range("S19").Value = 700
$ws.Range("T19").Value = 14
$ws.Range("D20").Value = 44511
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 24
$ws.Range("N20").Value = 140000
$ws.Range("O20").Value = 150000
$ws.Range("P20").Value = 145000
$ws.Range("S20").Value = 414
$ws.Range("D21").Value = 44193
$ws.Range("M21").Value = 8
$ws.Range("N21").Value = 150000
$ws.Range("O21").Value = 150000
$ws.Range("P21").Value = 150000
$ws.Range("S21").Value = 429
$ws.Range("D22").Value = 44167
$ws.Range("M22").Value = 140
$ws.Range("N22").Value = 9800
$ws.Range("O22").Value = 9800
$ws.Range("P22").Value = 9800
$ws.Range("Q22").Value = '$/caja 14 kilos empedrada'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 700
$ws.Range("T22").Value = 14
$ws.Range("D24").Value = 44189
$ws.Range("M24").Value = 16
$ws.Range("N24").Value = 150000
$ws.Range("P24").Value = 150000
$ws.Range("R24").Value = 'Provincia de Limarí'
$ws.Range("S24").Value = 429
